$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.81"
$ws.Range("E2").Value = "'-0.86%"
$ws.Range("D3").Value = "'35.39"
$ws.Range("E3").Value = "'-2.80%"
$ws.Range("D4").Value = "'5.083"
$ws.Range("E4").Value = "'0.30%"
$ws.Range("D5").Value = "'0.07989"
$ws.Range("E5").Value = "'0.73%"
$ws.Range("D6").Value = "'1.965"
$ws.Range("E6").Value = "'-10.91%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'7.841"
$ws.Range("E7").Value = "'-2.32%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.896"
$ws.Range("E8").Value = "'10.04%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9250"
$ws.Range("E9").Value = "'-0.60%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1092"
$ws.Range("E10").Value = "'10.80%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1892"
$ws.Range("E11").Value = "'0.87%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09509"
$ws.Range("E12").Value = "'4.63%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03662"
$ws.Range("E13").Value = "'0.28%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09918"
$ws.Range("E14").Value = "'-0.09%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001409"
$ws.Range("E15").Value = "'-2.29%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005790"
$ws.Range("E16").Value = "'2.89%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.456"
$ws.Range("E17").Value = "'-0.68%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.121"
$ws.Range("E18").Value = "'-0.98%"
$ws.Range("E19").Value = "'1.34%"
$ws.Range("E20").Value = "'-2.73%"
$ws.Range("D21").Value = "'5.115"
$ws.Range("E21").Value = "'-0.22%"
$ws.Range("E22").Value = "'0.16%"
$ws.Range("D23").Value = "'0.04523"
$ws.Range("E23").Value = "'-0.98%"
$ws.Range("D24").Value = "'0.001228"
$ws.Range("E24").Value = "'-1.01%"
$ws.Range("D25").Value = "'0.004691"
$ws.Range("E25").Value = "'-1.76%"
$ws.Range("D26").Value = "'0.0001258"
$ws.Range("E26").Value = "'-3.30%"
$ws.Range("D27").Value = "'0.0004455"
$ws.Range("E27").Value = "'-5.97%"
$ws.Range("D39").Value = "'0.01895"
$ws.Range("E39").Value = "'-3.98%"
$ws.Range("D40").Value = "'0.04740"
$ws.Range("E40").Value = "'-3.81%"
$ws.Range("D41").Value = "'0.007591"
$ws.Range("E41").Value = "'-2.75%"
$ws.Range("D42").Value = "'0.009656"
$ws.Range("E42").Value = "'23.62%"
$ws.Range("D43").Value = "'0.1341"
$ws.Range("E43").Value = "'-3.82%"
$ws.Range("D44").Value = "'0.002125"
$ws.Range("E44").Value = "'0.60%"
$ws.Range("D45").Value = "'0.01137"
$ws.Range("E45").Value = "'1.43%"
$ws.Range("D46").Value = "'0.00006258"
$ws.Range("E46").Value = "'0.59%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.12%"
$ws.Range("D48").Value = "'64.38"
$ws.Range("E48").Value = "'23.97%"
$ws.Range("D49").Value = "'0.001302"
$ws.Range("E49").Value = "'-27.77%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.12%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.12%"
